$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 20)
$ws.Range("D2").Value2 = 44893
$ws.Range("J2").Value2 = 900
$ws.Range("K2").Value2 = 13000
$ws.Range("L2").Value2 = 14000
$ws.Range("M2").Value2 = 13444
$ws.Range("P2").Value2 = 1034

# Row 3 (was row 15)
$ws.Range("D3").Value2 = 44984
$ws.Range("J3").Value2 = 400
$ws.Range("K3").Value2 = 16000
$ws.Range("L3").Value2 = 17000
$ws.Range("M3").Value2 = 16500
$ws.Range("P3").Value2 = 1269

# Row 4 (was row 17)
$ws.Range("D4").Value2 = 45141
$ws.Range("J4").Value2 = 400
$ws.Range("K4").Value2 = 16000
$ws.Range("L4").Value2 = 17000
$ws.Range("M4").Value2 = 16550
$ws.Range("N4").Value2 = "`$/caja 13 kilos"
$ws.Range("P4").Value2 = 1273
$ws.Range("Q4").Value2 = 13

# Row 5 (was row 18)
$ws.Range("D5").Value2 = 45154
$ws.Range("J5").Value2 = 250
$ws.Range("K5").Value2 = 17000
$ws.Range("L5").Value2 = 18000
$ws.Range("M5").Value2 = 17500
$ws.Range("P5").Value2 = 1346

# Row 6 (was row 19)
$ws.Range("D6").Value2 = 45155
$ws.Range("J6").Value2 = 300
$ws.Range("K6").Value2 = 16000
$ws.Range("L6").Value2 = 17000
$ws.Range("M6").Value2 = 16500
$ws.Range("P6").Value2 = 1269

# Row 7 (was row 10)
$ws.Range("D7").Value2 = 44469
$ws.Range("J7").Value2 = 140
$ws.Range("K7").Value2 = 13000
$ws.Range("L7").Value2 = 14000
$ws.Range("M7").Value2 = 13500
$ws.Range("P7").Value2 = 1038

# Row 8 (was row 2)
$ws.Range("D8").Value2 = 44406
$ws.Range("J8").Value2 = 160
$ws.Range("K8").Value2 = 17000
$ws.Range("L8").Value2 = 18000
$ws.Range("M8").Value2 = 17500
$ws.Range("P8").Value2 = 1346

# Row 9 (was row 7)
$ws.Range("D9").Value2 = 44592
$ws.Range("J9").Value2 = 120
$ws.Range("K9").Value2 = 12000
$ws.Range("L9").Value2 = 13000
$ws.Range("M9").Value2 = 12500
$ws.Range("P9").Value2 = 962

# Row 10 (was row 26)
$ws.Range("D10").Value2 = 44320
$ws.Range("J10").Value2 = 160
$ws.Range("K10").Value2 = 19000
$ws.Range("L10").Value2 = 20000
$ws.Range("M10").Value2 = 19500
$ws.Range("P10").Value2 = 1500

# Row 11 (was row 31)
$ws.Range("D11").Value2 = 45028
$ws.Range("J11").Value2 = 300
$ws.Range("K11").Value2 = 14000
$ws.Range("L11").Value2 = 15000
$ws.Range("M11").Value2 = 14500
$ws.Range("P11").Value2 = 1115

# Row 12 (was row 14)
$ws.Range("D12").Value2 = 45049
$ws.Range("K12").Value2 = 13000
$ws.Range("L12").Value2 = 14000
$ws.Range("M12").Value2 = 13500
$ws.Range("P12").Value2 = 1038

# Row 13 (was row 22)
$ws.Range("D13").Value2 = 45096
$ws.Range("K13").Value2 = 14000
$ws.Range("L13").Value2 = 15000
$ws.Range("M13").Value2 = 14600
$ws.Range("P13").Value2 = 1123

# Row 14 (was row 4)
$ws.Range("D14").Value2 = 44972
$ws.Range("J14").Value2 = 350
$ws.Range("K14").Value2 = 17000
$ws.Range("L14").Value2 = 18000
$ws.Range("M14").Value2 = 17429
$ws.Range("N14").Value2 = "`$/caja 15 kilos"
$ws.Range("P14").Value2 = 1162
$ws.Range("Q14").Value2 = 15

# Row 15 (was row 6)
$ws.Range("D15").Value2 = 45092
$ws.Range("J15").Value2 = 600
$ws.Range("K15").Value2 = 13000
$ws.Range("L15").Value2 = 14000
$ws.Range("M15").Value2 = 13500
$ws.Range("P15").Value2 = 1038

# Row 16 (was row 30)
$ws.Range("D16").Value2 = 44616
$ws.Range("J16").Value2 = 120
$ws.Range("K16").Value2 = 19000
$ws.Range("L16").Value2 = 20000
$ws.Range("M16").Value2 = 19500
$ws.Range("P16").Value2 = 1500

# Row 17 (was row 32)
$ws.Range("D17").Value2 = 44914
$ws.Range("J17").Value2 = 100
$ws.Range("K17").Value2 = 14000
$ws.Range("L17").Value2 = 15000
$ws.Range("M17").Value2 = 14400
$ws.Range("P17").Value2 = 1108

# Row 18 (was row 12)
$ws.Range("D18").Value2 = 44890
$ws.Range("J18").Value2 = 300
$ws.Range("K18").Value2 = 14000
$ws.Range("L18").Value2 = 15000
$ws.Range("M18").Value2 = 14500
$ws.Range("P18").Value2 = 1115

# Row 19 (was row 8)
$ws.Range("D19").Value2 = 44910
$ws.Range("J19").Value2 = 50
$ws.Range("K19").Value2 = 14000
$ws.Range("L19").Value2 = 15000
$ws.Range("M19").Value2 = 14500
$ws.Range("P19").Value2 = 1115

# Row 20 (was row 13)
$ws.Range("D20").Value2 = 44988
$ws.Range("J20").Value2 = 750
$ws.Range("K20").Value2 = 17000
$ws.Range("L20").Value2 = 18000
$ws.Range("M20").Value2 = 17400
$ws.Range("P20").Value2 = 1338

# Row 21 (was row 29)
$ws.Range("D21").Value2 = 44918
$ws.Range("I21").Value2 = "Segunda"
$ws.Range("J21").Value2 = 200
$ws.Range("K21").Value2 = 12000
$ws.Range("L21").Value2 = 13000
$ws.Range("M21").Value2 = 12750
$ws.Range("P21").Value2 = 981

# Row 22 (was row 16)
$ws.Range("D22").Value2 = 44397
$ws.Range("J22").Value2 = 140
$ws.Range("K22").Value2 = 12500
$ws.Range("L22").Value2 = 13000
$ws.Range("M22").Value2 = 12750
$ws.Range("P22").Value2 = 981

# Row 24 (was row 27)
$ws.Range("D24").Value2 = 44943
$ws.Range("I24").Value2 = "Segunda"
$ws.Range("J24").Value2 = 350
$ws.Range("K24").Value2 = 14000
$ws.Range("L24").Value2 = 15000
$ws.Range("M24").Value2 = 14429
$ws.Range("P24").Value2 = 1110

# Row 25 (was row 28)
$ws.Range("D25").Value2 = 45140
$ws.Range("K25").Value2 = 16000
$ws.Range("L25").Value2 = 17000
$ws.Range("M25").Value2 = 16500
$ws.Range("P25").Value2 = 1269

# Row 26 (was row 5)
$ws.Range("D26").Value2 = 44159
$ws.Range("J26").Value2 = 100
$ws.Range("K26").Value2 = 23000
$ws.Range("L26").Value2 = 24000
$ws.Range("M26").Value2 = 23500
$ws.Range("P26").Value2 = 1808

# Row 27 (was row 9)
$ws.Range("D27").Value2 = 44580
$ws.Range("I27").Value2 = "Primera"
$ws.Range("J27").Value2 = 160
$ws.Range("K27").Value2 = 11000
$ws.Range("L27").Value2 = 12000
$ws.Range("M27").Value2 = 11500
$ws.Range("P27").Value2 = 885

# Row 28 (was row 21)
$ws.Range("D28").Value2 = 44832
$ws.Range("J28").Value2 = 100
$ws.Range("K28").Value2 = 13000
$ws.Range("L28").Value2 = 14000
$ws.Range("M28").Value2 = 13500
$ws.Range("P28").Value2 = 1038

# Row 29 (was row 11)
$ws.Range("D29").Value2 = 44389
$ws.Range("I29").Value2 = "Primera"
$ws.Range("J29").Value2 = 120
$ws.Range("M29").Value2 = 12500
$ws.Range("P29").Value2 = 962

# Row 30 (was row 24)
$ws.Range("D30").Value2 = 44855
$ws.Range("J30").Value2 = 500
$ws.Range("K30").Value2 = 10000
$ws.Range("L30").Value2 = 10000
$ws.Range("M30").Value2 = 10000
$ws.Range("P30").Value2 = 769

# Row 31 (was row 25)
$ws.Range("D31").Value2 = 44764
$ws.Range("J31").Value2 = 200
$ws.Range("K31").Value2 = 12000
$ws.Range("L31").Value2 = 13000
$ws.Range("M31").Value2 = 12500
$ws.Range("P31").Value2 = 962

# Row 32 (was row 3)
$ws.Range("D32").Value2 = 44379
$ws.Range("J32").Value2 = 120
$ws.Range("K32").Value2 = 12000
$ws.Range("L32").Value2 = 13000
$ws.Range("M32").Value2 = 12667
$ws.Range("P32").Value2 = 974
